$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column Q (previous year, 2020) into new column R (2021) for rows 3-34
$ws.Range("Q3:Q34").Copy()
$ws.Range("R3:R34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new 2021 data values for column R
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 109
$ws.Range("R6").Value = 74
$ws.Range("R7").Value = 35
$ws.Range("R8").Value = 36
$ws.Range("R9").Value = 35
$ws.Range("R10").Value = 1
$ws.Range("R11").Value = 15
$ws.Range("R12").Value = 8
$ws.Range("R13").Value = 7
$ws.Range("R14").Value = 12
$ws.Range("R15").Value = 7
$ws.Range("R16").Value = 5
$ws.Range("R17").Value = "-"
$ws.Range("R18").Value = "-"
$ws.Range("R19").Value = "-"
$ws.Range("R20").Value = 17
$ws.Range("R21").Value = 8
$ws.Range("R22").Value = 9
$ws.Range("R23").Value = 9
$ws.Range("R24").Value = 7
$ws.Range("R25").Value = 2
$ws.Range("R26").Value = 20
$ws.Range("R27").Value = 9
$ws.Range("R28").Value = 11
$ws.Range("R29").Value = "-"
$ws.Range("R30").Value = "-"
$ws.Range("R31").Value = "-"
$ws.Range("R32").Value = "-"
$ws.Range("R33").Value = "-"
$ws.Range("R34").Value = "-"

# Update the selected cell to reflect the new last column, matching the saved view state
$ws.Range("R1").Select() | Out-Null
